# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp update
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 19:52"

# Row 4: Estados Unidos - updated totals
$ws.Range("B4").Value = 418044
$ws.Range("C4").Value = 17709
$ws.Range("D4").Value = 22184
$ws.Range("E4").Value = 381646
$ws.Range("F4").Value = 9224
$ws.Range("G4").Value = 1373
$ws.Range("H4").Value = 14214

# Rows 7 & 8: Francia overtakes Alemania in the ranking, so the two
# countries (and their data) swap places.
# Row 7 becomes Francia with fresh updated figures.
$ws.Range("A7").Value = "Francia"
$ws.Range("B7").Value = 112950
$ws.Range("C7").Value = 3881
$ws.Range("D7").Value = 21254
$ws.Range("E7").Value = 80827
$ws.Range("F7").Value = 7148
$ws.Range("G7").Value = 541
$ws.Range("H7").Value = 10869

# Row 8 becomes Alemania, carrying what used to be row 7's figures.
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 109702
$ws.Range("C8").Value = 2039
$ws.Range("D8").Value = 36081
$ws.Range("E8").Value = 71516
$ws.Range("F8").Value = 4895
$ws.Range("G8").Value = 89
$ws.Range("H8").Value = 2105

# Row 60: Marruecos
$ws.Range("B60").Value = 1275
$ws.Range("C60").Value = 91
$ws.Range("E60").Value = 1085
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 93

# Row 100: Malta
$ws.Range("E100").Value = 293
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 1

# Row 186: Groenlandia
$ws.Range("D186").Value = 11
$ws.Range("E186").Value = 0
